$wb = $excel.ActiveWorkbook

# New row (row 40) data for each of the 4 worksheets (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2)
# Columns: A=time, B=totalLen, C=ID, D=actualLen, E=checksum, F=totalLen_DEC, G=ID_DEC, H=actualLen_DEC, I=checksum_DEC

$rowsData = @(
    @{ Sheet = 1; A = 45826.4603587963; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x7C"; E = "0x07"; F = 400; G = "5.68631262647113e+23"; H = 380; I = 7 },
    @{ Sheet = 2; A = 45826.4603587963; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x6C"; E = "0x19"; F = 380; G = "5.68432987514711e+23"; H = 364; I = 25 },
    @{ Sheet = 3; A = 45826.4603587963; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x6B"; E = "0x15"; F = 110; G = "5.68631262647113e+23"; H = 107; I = 15 },
    @{ Sheet = 4; A = 45826.4603587963; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x7F"; E = "0x9"; F = 130; G = "5.68631262647113e+23"; H = 127; I = 9 }
)

foreach ($rowData in $rowsData) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)
    $newRow = 40

    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($newRow, 1).Value = $rowData.A
    $ws.Cells.Item($newRow, 2).Value = $rowData.B
    $ws.Cells.Item($newRow, 3).Value = $rowData.C
    $ws.Cells.Item($newRow, 4).Value = $rowData.D
    $ws.Cells.Item($newRow, 5).Value = $rowData.E
    $ws.Cells.Item($newRow, 6).Value = $rowData.F
    $ws.Cells.Item($newRow, 7).Value = [double]$rowData.G
    $ws.Cells.Item($newRow, 8).Value = $rowData.H
    $ws.Cells.Item($newRow, 9).Value = $rowData.I
}
